# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.261.09"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.449.91"
$ws.Range("E3").Value = "  +2.91%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9427"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.029"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06517"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9979"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.353"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.070"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "1.444.35"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9573"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05693"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("E21").Value = "  -4.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.47%  "
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.231"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "20.289.99"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.084"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").Value = "1.596.03"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.952"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.810"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.04%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7866"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07712"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.477"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05643"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.656"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.114"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9483"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1852"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.373"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5219"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.473"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5104"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06378"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9831"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.05%  "
